# Femacal de La Calera - Chirimoya: add a new weekly price report block.
#
# The source publishes one 3-row block (Especial/Primera/Segunda) per market
# day, newest block first (just under the header + the most recent block
# already present). This edit inserts a brand-new block for 2022-11-10
# (serial 44875) just above the existing history, pushing every existing
# data row down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 193 (format is copied from the
# row above, same as Excel's native "Insert Copied Cells"/"Insert Sheet
# Rows" behaviour), shifting the old rows 193:263 down to 196:266.
$ws.Rows("193:195").Insert()

# Row 193: Especial
$ws.Range("A193").Value = 3
$ws.Range("B193").Value = "Femacal de La Calera"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44875
$ws.Range("E193").Value = 5
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100107
$ws.Range("H193").Value = "Otros"
$ws.Range("I193").Value = 100107002
$ws.Range("J193").Value = "Chirimoya"
$ws.Range("K193").Value = "Cultivar IV Región"
$ws.Range("L193").Value = "Especial"
$ws.Range("M193").Value = 60
$ws.Range("N193").Value = 28000
$ws.Range("O193").Value = 28000
$ws.Range("P193").Value = 28000
$ws.Range("Q193").Value = "$/bandeja 10 kilos"
$ws.Range("R193").Value = "Provincia del Elquí"
$ws.Range("S193").Value = 2800
$ws.Range("T193").Value = 10

# Row 194: Primera
$ws.Range("A194").Value = 3
$ws.Range("B194").Value = "Femacal de La Calera"
$ws.Range("C194").Value = "Coquimbo"
$ws.Range("D194").Value = 44875
$ws.Range("E194").Value = 5
$ws.Range("F194").Value = "Fruta"
$ws.Range("G194").Value = 100107
$ws.Range("H194").Value = "Otros"
$ws.Range("I194").Value = 100107002
$ws.Range("J194").Value = "Chirimoya"
$ws.Range("K194").Value = "Cultivar IV Región"
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 56
$ws.Range("N194").Value = 25000
$ws.Range("O194").Value = 25000
$ws.Range("P194").Value = 25000
$ws.Range("Q194").Value = "$/bandeja 10 kilos"
$ws.Range("R194").Value = "Provincia del Elquí"
$ws.Range("S194").Value = 2500
$ws.Range("T194").Value = 10

# Row 195: Segunda
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44875
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100107
$ws.Range("H195").Value = "Otros"
$ws.Range("I195").Value = 100107002
$ws.Range("J195").Value = "Chirimoya"
$ws.Range("K195").Value = "Cultivar IV Región"
$ws.Range("L195").Value = "Segunda"
$ws.Range("M195").Value = 57
$ws.Range("N195").Value = 22000
$ws.Range("O195").Value = 22000
$ws.Range("P195").Value = 22000
$ws.Range("Q195").Value = "$/bandeja 10 kilos"
$ws.Range("R195").Value = "Provincia del Elquí"
$ws.Range("S195").Value = 2200
$ws.Range("T195").Value = 10
